$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 288896462.044739604
$ws.Range("C2").Value = 341886003.8149094582
$ws.Range("D2").Value = 394875545.5850800276
$ws.Range("E2").Value = 447865087.3552497625
$ws.Range("F2").Value = 500854629.1254184842
$ws.Range("B3").Value = 683889683.8933423758
$ws.Range("C3").Value = 736879225.6635122299
$ws.Range("D3").Value = 789868767.4336826801
$ws.Range("E3").Value = 842858309.2038526535
$ws.Range("F3").Value = 895847850.9740213156
$ws.Range("B4").Value = 1474229291.0230967999
$ws.Range("C4").Value = 1527218832.7932667732
$ws.Range("D4").Value = 1580208374.5634374619
$ws.Range("E4").Value = 1633197916.3336069584
$ws.Range("F4").Value = 1686187458.1037759781
$ws.Range("B5").Value = 2423258707.2723855972
$ws.Range("C5").Value = 2476248249.042555809
$ws.Range("D5").Value = 2529237790.8127264977
$ws.Range("E5").Value = 2582227332.5828962326
$ws.Range("F5").Value = 2635216874.3530650139
